$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Date fix: "1/9/2015" -> "1/10/2015", leaving the "Date:" run untouched,
#    and leaving the user's edit cursor (the _GoBack bookmark) right after
#    the newly-typed "1/10", splitting the date run into "  1/10" / "/2015".
# ---------------------------------------------------------------------------

# Temporary barrier right after "Date:" so the engine's run-merge pass
# doesn't fold the "Date:" run into the date-value run while we edit it.
$dateLabel = $d.Content
$dateLabel.Find.Execute("Date:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateBarrierPos = $dateLabel.End
$d.Bookmarks.Add("ZZZ_DateBarrier", $d.Range($dateBarrierPos, $dateBarrierPos))

$dateRng = $d.Content
$dateRng.Find.Execute("1/9", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateRng.Text = "1/10"

$d.Bookmarks("ZZZ_DateBarrier").Delete()

# Relocate _GoBack to sit right after the newly-typed "1/10" (before "/2015"),
# matching where Word leaves the last-edit marker.
$afterTyped = $d.Content
$afterTyped.Find.Execute("1/10", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPos = $afterTyped.End
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))

# ---------------------------------------------------------------------------
# 2) Project-sketch sentence: accept the grammar check on "looks like" by
#    merging " an image of what your project sketch ", "looks" and " like."
#    into a single run (this also clears the w:proofErr gramStart/gramEnd
#    markers Word had inserted around "looks").
# ---------------------------------------------------------------------------

# Barriers keep the untouched neighboring runs ("insert" and " Make sure to ")
# from being swept into the merge.
$beforeImage = $d.Content
$beforeImage.Find.Execute("insert", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$barrier1Pos = $beforeImage.End
$d.Bookmarks.Add("ZZZ_Barrier1", $d.Range($barrier1Pos, $barrier1Pos))

$afterLike = $d.Content
$afterLike.Find.Execute(" like.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$barrier2Pos = $afterLike.End
$d.Bookmarks.Add("ZZZ_Barrier2", $d.Range($barrier2Pos, $barrier2Pos))

$sketchRng = $d.Content
$sketchRng.Find.Execute(" an image of what your project sketch looks like.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sketchRng.Text = " an image of what your project sketch looks like.#"

$sketchRng2 = $d.Content
$sketchRng2.Find.Execute(" an image of what your project sketch looks like.#", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sketchRng2.Text = " an image of what your project sketch looks like."

$d.Bookmarks("ZZZ_Barrier1").Delete()
$d.Bookmarks("ZZZ_Barrier2").Delete()
